$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I8").Value = 0.468739117536998
$ws.Range("J8").Value = 0.2150495036779461

$ws.Range("H9").Value = 0.5099036351493167
$ws.Range("I9").Value = 0.24

$ws.Range("G10").Value = 0.5604363747513331
$ws.Range("H10").Value = 0.3087982760018804

$ws.Range("F11").Value = 0.5999036351493168
$ws.Range("G11").Value = 0.32

$ws.Range("E12").Value = 0.6299036351493167
$ws.Range("F12").Value = 0.4476495795507702

$ws.Range("D13").Value = 0.3603773643037867
$ws.Range("E13").Value = 0.1088966743764388

$ws.Range("C14").Value = 0.4107440146302961
$ws.Range("D14").Value = 0.1461563307127136

$ws.Range("B15").Value = 0.25708246933236
$ws.Range("C15").Value = 0.09547648014918764

$ws.Range("B16").Value = 0.0959495356205764
